$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New s_val data (regenerated to filter save games)
$data = @{
    2 = @(0.01253208636536152, 0.04103571897497393, 0.7210945179870265, 13.86384647080068, 14.63850879412805)
    3 = @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1.104883657715537)
    4 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    5 = @(1.445647641019636, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 2.433531715253719)
    6 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    7 = @(0.1169995834814548, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.426980108624251)
    8 = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 2.213936997104367)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
